$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell D1 = "success", reusing the same header style as C1
# (copy C1's formatting into D1 first, then overwrite the value so the
# cell keeps style index 1 instead of Excel minting a brand-new xf).
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = "success"

# Fill D2:D38 with the textual "0"/"1" success flag (row 13 -> "1", rest -> "0").
# Writing the literal string "0" via .Value gets auto-coerced to a number by
# COM, so instead we materialize it through a TEXT() formula and then paste
# the computed value back over itself (values-only), which preserves the
# Text cell type without touching cell formatting/styles.
for ($r = 2; $r -le 38; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($r -eq 13) {
        $cell.Formula = "=TEXT(1,""0"")"
    } else {
        $cell.Formula = "=TEXT(0,""0"")"
    }
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}
